# Update header labels on existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new "PO Forecast" sheet after the existing sheets
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$wsForecast.Cells.Item(2, 1).Value = 45228.99999999999
$wsForecast.Cells.Item(2, 2).Value = 150
$wsForecast.Cells.Item(2, 3).Value = -38.82554119131935
$wsForecast.Cells.Item(2, 4).Value = 322.6024916803933
$wsForecast.Cells.Item(3, 1).Value = 45242.99999999999
$wsForecast.Cells.Item(3, 2).Value = 149
$wsForecast.Cells.Item(3, 3).Value = -35.9412143418525
$wsForecast.Cells.Item(3, 4).Value = 322.6365226241792
$wsForecast.Cells.Item(4, 1).Value = 45249.99999999999
$wsForecast.Cells.Item(4, 2).Value = 148
$wsForecast.Cells.Item(4, 3).Value = -34.88681792021994
$wsForecast.Cells.Item(4, 4).Value = 323.2388326719374
$wsForecast.Cells.Item(5, 1).Value = 45256.99999999999
$wsForecast.Cells.Item(5, 2).Value = 147
$wsForecast.Cells.Item(5, 3).Value = -43.70693769658835
$wsForecast.Cells.Item(5, 4).Value = 333.1334006695303
$wsForecast.Cells.Item(6, 1).Value = 45263.99999999999
$wsForecast.Cells.Item(6, 2).Value = 146
$wsForecast.Cells.Item(6, 3).Value = -32.09783240312655
$wsForecast.Cells.Item(6, 4).Value = 326.9300110043043
$wsForecast.Cells.Item(7, 1).Value = 45298.99999999999
$wsForecast.Cells.Item(7, 2).Value = 142
$wsForecast.Cells.Item(7, 3).Value = -47.12031871403227
$wsForecast.Cells.Item(7, 4).Value = 323.4814669115495
$wsForecast.Cells.Item(8, 1).Value = 45305.99999999999
$wsForecast.Cells.Item(8, 2).Value = 142
$wsForecast.Cells.Item(8, 3).Value = -48.66740581112752
$wsForecast.Cells.Item(8, 4).Value = 314.7568334265616
$wsForecast.Cells.Item(9, 1).Value = 45312.99999999999
$wsForecast.Cells.Item(9, 2).Value = 141
$wsForecast.Cells.Item(9, 3).Value = -35.7472086156849
$wsForecast.Cells.Item(9, 4).Value = 329.1545284556537
$wsForecast.Cells.Item(10, 1).Value = 45319.99999999999
$wsForecast.Cells.Item(10, 2).Value = 140
$wsForecast.Cells.Item(10, 3).Value = -44.79251488286394
$wsForecast.Cells.Item(10, 4).Value = 321.3068678917327
$wsForecast.Cells.Item(11, 1).Value = 45326.99999999999
$wsForecast.Cells.Item(11, 2).Value = 139
$wsForecast.Cells.Item(11, 3).Value = -37.60001195409031
$wsForecast.Cells.Item(11, 4).Value = 330.9467085444097
$wsForecast.Cells.Item(12, 1).Value = 45333.99999999999
$wsForecast.Cells.Item(12, 2).Value = 138
$wsForecast.Cells.Item(12, 3).Value = -51.50067434717143
$wsForecast.Cells.Item(12, 4).Value = 320.8829755452511
$wsForecast.Cells.Item(13, 1).Value = 45340.99999999999
$wsForecast.Cells.Item(13, 2).Value = 138
$wsForecast.Cells.Item(13, 3).Value = -65.00860894387027
$wsForecast.Cells.Item(13, 4).Value = 319.0950799925637
$wsForecast.Cells.Item(14, 1).Value = 45347.99999999999
$wsForecast.Cells.Item(14, 2).Value = 137
$wsForecast.Cells.Item(14, 3).Value = -43.13014793157149
$wsForecast.Cells.Item(14, 4).Value = 313.6841327256521
$wsForecast.Cells.Item(15, 1).Value = 45361.99999999999
$wsForecast.Cells.Item(15, 2).Value = 135
$wsForecast.Cells.Item(15, 3).Value = -51.02535204314236
$wsForecast.Cells.Item(15, 4).Value = 304.4254293743263
$wsForecast.Cells.Item(16, 1).Value = 45368.99999999999
$wsForecast.Cells.Item(16, 2).Value = 134
$wsForecast.Cells.Item(16, 3).Value = -45.15420818849563
$wsForecast.Cells.Item(16, 4).Value = 315.3619517513966
$wsForecast.Cells.Item(17, 1).Value = 45375.99999999999
$wsForecast.Cells.Item(17, 2).Value = 134
$wsForecast.Cells.Item(17, 3).Value = -51.16715861797581
$wsForecast.Cells.Item(17, 4).Value = 310.12980293002
$wsForecast.Cells.Item(18, 1).Value = 45382.99999999999
$wsForecast.Cells.Item(18, 2).Value = 133
$wsForecast.Cells.Item(18, 3).Value = -54.20172628090079
$wsForecast.Cells.Item(18, 4).Value = 315.11978533226
$wsForecast.Cells.Item(19, 1).Value = 45389.99999999999
$wsForecast.Cells.Item(19, 2).Value = 132
$wsForecast.Cells.Item(19, 3).Value = -47.91283506349123
$wsForecast.Cells.Item(19, 4).Value = 302.2986735327142
$wsForecast.Cells.Item(20, 1).Value = 45396.99999999999
$wsForecast.Cells.Item(20, 2).Value = 131
$wsForecast.Cells.Item(20, 3).Value = -50.69320666564465
$wsForecast.Cells.Item(20, 4).Value = 303.4579801501862
$wsForecast.Cells.Item(21, 1).Value = 45403.99999999999
$wsForecast.Cells.Item(21, 2).Value = 130
$wsForecast.Cells.Item(21, 3).Value = -42.1852297309227
$wsForecast.Cells.Item(21, 4).Value = 322.3585691520506
$wsForecast.Cells.Item(22, 1).Value = 45410.99999999999
$wsForecast.Cells.Item(22, 2).Value = 130
$wsForecast.Cells.Item(22, 3).Value = -61.68792945779048
$wsForecast.Cells.Item(22, 4).Value = 303.2413738465746
$wsForecast.Cells.Item(23, 1).Value = 45417.99999999999
$wsForecast.Cells.Item(23, 2).Value = 129
$wsForecast.Cells.Item(23, 3).Value = -56.53904289857172
$wsForecast.Cells.Item(23, 4).Value = 313.2495093178495
$wsForecast.Cells.Item(24, 1).Value = 45424.99999999999
$wsForecast.Cells.Item(24, 2).Value = 128
$wsForecast.Cells.Item(24, 3).Value = -45.1982524069681
$wsForecast.Cells.Item(24, 4).Value = 313.4916527362117
$wsForecast.Cells.Item(25, 1).Value = 45431.99999999999
$wsForecast.Cells.Item(25, 2).Value = 127
$wsForecast.Cells.Item(25, 3).Value = -54.00942284398861
$wsForecast.Cells.Item(25, 4).Value = 307.7284699362143
$wsForecast.Cells.Item(26, 1).Value = 45445.99999999999
$wsForecast.Cells.Item(26, 2).Value = 126
$wsForecast.Cells.Item(26, 3).Value = -46.30245902683105
$wsForecast.Cells.Item(26, 4).Value = 295.6438240247406
$wsForecast.Cells.Item(27, 1).Value = 45459.99999999999
$wsForecast.Cells.Item(27, 2).Value = 124
$wsForecast.Cells.Item(27, 3).Value = -58.5963628770498
$wsForecast.Cells.Item(27, 4).Value = 298.2470083850637
$wsForecast.Cells.Item(28, 1).Value = 45466.99999999999
$wsForecast.Cells.Item(28, 2).Value = 123
$wsForecast.Cells.Item(28, 3).Value = -43.54875628988244
$wsForecast.Cells.Item(28, 4).Value = 303.8266063740492
$wsForecast.Cells.Item(29, 1).Value = 45480.99999999999
$wsForecast.Cells.Item(29, 2).Value = 122
$wsForecast.Cells.Item(29, 3).Value = -50.67383667816667
$wsForecast.Cells.Item(29, 4).Value = 304.178918878507
$wsForecast.Cells.Item(30, 1).Value = 45487.99999999999
$wsForecast.Cells.Item(30, 2).Value = 121
$wsForecast.Cells.Item(30, 3).Value = -54.38261538113193
$wsForecast.Cells.Item(30, 4).Value = 300.3961905244024
$wsForecast.Cells.Item(31, 1).Value = 45501.99999999999
$wsForecast.Cells.Item(31, 2).Value = 119
$wsForecast.Cells.Item(31, 3).Value = -69.23966595451357
$wsForecast.Cells.Item(31, 4).Value = 295.4000547114027
$wsForecast.Cells.Item(32, 1).Value = 45508.99999999999
$wsForecast.Cells.Item(32, 2).Value = 118
$wsForecast.Cells.Item(32, 3).Value = -61.71397370747485
$wsForecast.Cells.Item(32, 4).Value = 298.8409581774545
$wsForecast.Cells.Item(33, 1).Value = 45543.99999999999
$wsForecast.Cells.Item(33, 2).Value = 114
$wsForecast.Cells.Item(33, 3).Value = -64.67770759286327
$wsForecast.Cells.Item(33, 4).Value = 310.8444855401764
$wsForecast.Cells.Item(34, 1).Value = 45550.99999999999
$wsForecast.Cells.Item(34, 2).Value = 114
$wsForecast.Cells.Item(34, 3).Value = -71.82788784535126
$wsForecast.Cells.Item(34, 4).Value = 300.059294696376
$wsForecast.Cells.Item(35, 1).Value = 45557.99999999999
$wsForecast.Cells.Item(35, 2).Value = 113
$wsForecast.Cells.Item(35, 3).Value = -63.28733479282074
$wsForecast.Cells.Item(35, 4).Value = 295.613804473445
$wsForecast.Cells.Item(36, 1).Value = 45571.99999999999
$wsForecast.Cells.Item(36, 2).Value = 111
$wsForecast.Cells.Item(36, 3).Value = -69.23024657165493
$wsForecast.Cells.Item(36, 4).Value = 292.8645223227371
$wsForecast.Cells.Item(37, 1).Value = 45578.99999999999
$wsForecast.Cells.Item(37, 2).Value = 110
$wsForecast.Cells.Item(37, 3).Value = -60.7141428993279
$wsForecast.Cells.Item(37, 4).Value = 283.6492877782957
$wsForecast.Cells.Item(38, 1).Value = 45585.99999999999
$wsForecast.Cells.Item(38, 2).Value = 110
$wsForecast.Cells.Item(38, 3).Value = -69.30005789310522
$wsForecast.Cells.Item(38, 4).Value = 300.0450879158885
$wsForecast.Cells.Item(39, 1).Value = 45592.99999999999
$wsForecast.Cells.Item(39, 2).Value = 109
$wsForecast.Cells.Item(39, 3).Value = -77.09557908254828
$wsForecast.Cells.Item(39, 4).Value = 292.9178082218697
$wsForecast.Cells.Item(40, 1).Value = 45599.99999999999
$wsForecast.Cells.Item(40, 2).Value = 108
$wsForecast.Cells.Item(40, 3).Value = -56.9657984932936
$wsForecast.Cells.Item(40, 4).Value = 288.1835482787429
$wsForecast.Cells.Item(41, 1).Value = 45606.99999999999
$wsForecast.Cells.Item(41, 2).Value = 107
$wsForecast.Cells.Item(41, 3).Value = -75.65232272591676
$wsForecast.Cells.Item(41, 4).Value = 285.9452569766532
$wsForecast.Cells.Item(42, 1).Value = 45613.99999999999
$wsForecast.Cells.Item(42, 2).Value = 106
$wsForecast.Cells.Item(42, 3).Value = -85.17823296157727
$wsForecast.Cells.Item(42, 4).Value = 290.9837152460105
$wsForecast.Cells.Item(43, 1).Value = 45620.99999999999
$wsForecast.Cells.Item(43, 2).Value = 106
$wsForecast.Cells.Item(43, 3).Value = -74.85499864285508
$wsForecast.Cells.Item(43, 4).Value = 288.1558202843386
$wsForecast.Cells.Item(44, 1).Value = 45627.99999999999
$wsForecast.Cells.Item(44, 2).Value = 105
$wsForecast.Cells.Item(44, 3).Value = -74.78172307867477
$wsForecast.Cells.Item(44, 4).Value = 297.1357696136095
$wsForecast.Cells.Item(45, 1).Value = 45634.99999999999
$wsForecast.Cells.Item(45, 2).Value = 104
$wsForecast.Cells.Item(45, 3).Value = -76.38815415989356
$wsForecast.Cells.Item(45, 4).Value = 281.0383241742708
$wsForecast.Cells.Item(46, 1).Value = 45641.99999999999
$wsForecast.Cells.Item(46, 2).Value = 103
$wsForecast.Cells.Item(46, 3).Value = -76.03321130753217
$wsForecast.Cells.Item(46, 4).Value = 294.928678064139
$wsForecast.Cells.Item(47, 1).Value = 45648.99999999999
$wsForecast.Cells.Item(47, 2).Value = 102
$wsForecast.Cells.Item(47, 3).Value = -75.5264843240721
$wsForecast.Cells.Item(47, 4).Value = 292.7599449039957
$wsForecast.Cells.Item(48, 1).Value = 45655.99999999999
$wsForecast.Cells.Item(48, 2).Value = 102
$wsForecast.Cells.Item(48, 3).Value = -87.56761390696485
$wsForecast.Cells.Item(48, 4).Value = 285.4931419443886
$wsForecast.Cells.Item(49, 1).Value = 45662.99999999999
$wsForecast.Cells.Item(49, 2).Value = 101
$wsForecast.Cells.Item(49, 3).Value = -80.37593274885607
$wsForecast.Cells.Item(49, 4).Value = 287.0766046689053

# Date formatting for column A (data rows)
$wsForecast.Range("A2:A49").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Re-select A1 on the first sheet to mirror original active selection
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
